$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data to the sheet.
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "248.33") are forced to Text format first, then the
# explicit style is cleared again so the cell formatting matches the
# rest of the (unstyled) data cells.

$ws.Range("D2").Value = "34.779.68"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "1.878.06"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.691"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.51%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0740"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("D13").Value = "2.152.33"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.718"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "1.902.47"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "34.778.76"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "0.0₃0822"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.45%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("E30").Value = "  -5.57%  "
$ws.Range("D31").Value = "4.128.41"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0581"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.837"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -16.50%  "
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0657"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.285.74"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0772"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("E51").Value = "  -2.53%  "
